# Fruta / hortaliza, semanal
# Insert a new weekly record as row 41, pushing the existing rows
# (old 41-56) down to (42-57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 41..56 down to 42..57 by inserting a new row at 41.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with this week's record.
$ws.Cells.Item(41, 1).Value = 1
$ws.Cells.Item(41, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(41, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(41, 4).Value = 44559
$ws.Cells.Item(41, 5).Value = 15
$ws.Cells.Item(41, 6).Value = 100114001
$ws.Cells.Item(41, 7).Value = "Papa"
$ws.Cells.Item(41, 8).Value = "Asterix"
$ws.Cells.Item(41, 9).Value = "1a nueva(o)"
$ws.Cells.Item(41, 10).Value = 1000
$ws.Cells.Item(41, 11).Value = 13000
$ws.Cells.Item(41, 12).Value = 14000
$ws.Cells.Item(41, 13).Value = 13500
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(41, 16).Value = 540
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(41, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
